$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I10").Value = 3.5
$ws.Range("J10").Value = 3.25
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("AA10").Value = 23
$ws.Range("AJ10").Value = 41
$ws.Range("AP10").Value = 29
$ws.Range("AR10").Value = 81
